$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 271, shifting the existing rows 271:327
# down to 272:328 (dimension grows from A1:R327 to A1:R328).
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A271").Value = 4
$ws.Range("B271").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C271").Value = "Los Lagos"
$ws.Range("D271").Value = 44641
$ws.Range("E271").Value = 10
$ws.Range("F271").Value = 100114013
$ws.Range("G271").Value = "Zanahoria"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 250
$ws.Range("K271").Value = 12000
$ws.Range("L271").Value = 12000
$ws.Range("M271").Value = 12000
$ws.Range("N271").Value = "$/saco 20 kilos"
$ws.Range("O271").Value = "Chillán"
$ws.Range("P271").Value = 600
$ws.Range("Q271").Value = 20
$ws.Range("R271").Value = "Hortaliza"
